# Modificato C6 da 10nF a 100nF per possibile instabilità refererence
# C6 moves from the C2,C5,C6 (10nF) group into the C1,C3,C4 (100nF) group.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 3: C1,C3,C4  -> C1,C3,C4,C6 (qty 3 -> 4)
$ws.Range("B3").Value = 4
$ws.Range("C3").Value = "C1,C3,C4,C6"

# Row 4: C2,C5,C6 -> C2,C5 (qty 3 -> 2)
$ws.Range("B4").Value = 2
$ws.Range("C4").Value = "C2,C5"
